$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple userCount (column F) corrections ---
$ws.Range("F2").Value = 271
$ws.Range("F4").Value = 105
$ws.Range("F7").Value = 59
$ws.Range("F8").Value = 93
$ws.Range("F9").Value = 324
$ws.Range("F10").Value = 86
$ws.Range("F11").Value = 473
$ws.Range("F12").Value = 83
$ws.Range("F14").Value = 70
$ws.Range("F16").Value = 50

# --- Rows 5 and 6 swap identities (Franco Zanelatto <-> Jhamir D'Arrigo), ---
# --- with updated userCount figures for each player ---

# New row 5: Jhamir D'Arrigo
$ws.Range("A5").Value = "Jhamir D´Arrigo"
$ws.Range("B5").Value = "jhamir-d-arrigo"
$ws.Range("C5").Value = "J. D´Arrigo"
$ws.Range("D5").Value = "M"
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 98
$ws.Range("G5").Value = 1017760
$ws.Range("J5").Value = 49.542857142857
$ws.Range("K5").Value = 82.814285714286
$ws.Range("L5").Value = 21
$ws.Range("M5").Value = "Alianza Lima"

# New row 6: Franco Zanelatto
$ws.Range("A6").Value = "Franco Zanelatto"
$ws.Range("B6").Value = "franco-zanelatto"
$ws.Range("C6").Value = "F. Zanelatto"
$ws.Range("D6").Value = "M"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 329
$ws.Range("G6").Value = 973650
$ws.Range("J6").Value = 78.22499999999999
$ws.Range("K6").Value = 57.775
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = "Alianza Lima"
